# Auto-generated: update calibrated values for Uganda SE model input sheet
# (strategy_id-0) rows 2, 10, 11, 12 per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("strategy_id-0")

# Row 2
$ws.Range("J2").Value = 32.38718373
$ws.Range("K2").Value = 33.93561508
$ws.Range("L2").Value = 34.9982768
$ws.Range("M2").Value = 37.2045415
$ws.Range("N2").Value = 39.60004706
$ws.Range("O2").Value = 40.7687658
$ws.Range("P2").Value = 42.21058595
$ws.Range("Q2").Value = 44.14721689
$ws.Range("R2").Value = 46.50334192
$ws.Range("S2").Value = 53.89737328
$ws.Range("T2").Value = 62.46705563
$ws.Range("U2").Value = 72.39931747999999
$ws.Range("V2").Value = 83.91080896
$ws.Range("W2").Value = 97.25262758
$ws.Range("X2").Value = 112.7157954
$ws.Range("Y2").Value = 130.6376068
$ws.Range("Z2").Value = 140.0435145
$ws.Range("AA2").Value = 150.1266476
$ws.Range("AB2").Value = 160.9357662
$ws.Range("AC2").Value = 172.5231414
$ws.Range("AD2").Value = 184.9448075
$ws.Range("AE2").Value = 199.3705025
$ws.Range("AF2").Value = 214.9214017
$ws.Range("AG2").Value = 231.6852711
$ws.Range("AH2").Value = 249.7567222
$ws.Range("AI2").Value = 269.2377465
$ws.Range("AJ2").Value = 292.3921927
$ws.Range("AK2").Value = 317.5379213
$ws.Range("AL2").Value = 344.8461825
$ws.Range("AM2").Value = 374.5029542
$ws.Range("AN2").Value = 406.7102083
$ws.Range("AO2").Value = 441.6872862
$ws.Range("AP2").Value = 479.6723928
$ws.Range("AQ2").Value = 520.9242186
$ws.Range("AR2").Value = 565.7237014
$ws.Range("AS2").Value = 614.3759397

# Row 10
$ws.Range("J10").Value = 4.812222222
$ws.Range("K10").Value = 4.812222222
$ws.Range("L10").Value = 4.812222222
$ws.Range("M10").Value = 4.812222222
$ws.Range("N10").Value = 4.812222222
$ws.Range("O10").Value = 4.812222222
$ws.Range("P10").Value = 4.812222222
$ws.Range("Q10").Value = 4.812222222
$ws.Range("R10").Value = 4.812222222
$ws.Range("S10").Value = 4.812222222
$ws.Range("T10").Value = 4.812222222
$ws.Range("U10").Value = 4.812222222
$ws.Range("V10").Value = 4.812222222
$ws.Range("W10").Value = 4.812222222
$ws.Range("X10").Value = 4.812222222
$ws.Range("Y10").Value = 4.812222222
$ws.Range("Z10").Value = 4.812222222
$ws.Range("AA10").Value = 4.812222222
$ws.Range("AB10").Value = 4.812222222
$ws.Range("AC10").Value = 4.812222222
$ws.Range("AD10").Value = 4.812222222
$ws.Range("AE10").Value = 4.812222222
$ws.Range("AF10").Value = 4.812222222
$ws.Range("AG10").Value = 4.812222222
$ws.Range("AH10").Value = 4.812222222
$ws.Range("AI10").Value = 4.812222222
$ws.Range("AJ10").Value = 4.812222222
$ws.Range("AK10").Value = 4.812222222
$ws.Range("AL10").Value = 4.812222222
$ws.Range("AM10").Value = 4.812222222
$ws.Range("AN10").Value = 4.812222222
$ws.Range("AO10").Value = 4.812222222
$ws.Range("AP10").Value = 4.812222222
$ws.Range("AQ10").Value = 4.812222222
$ws.Range("AR10").Value = 4.812222222
$ws.Range("AS10").Value = 4.812222222

# Row 11
$ws.Range("J11").Value = 29792913.39
$ws.Range("K11").Value = 30678944.1
$ws.Range("L11").Value = 31617586.25
$ws.Range("M11").Value = 32570631.93
$ws.Range("N11").Value = 33485072.91
$ws.Range("P11").Value = 34862528.73
$ws.Range("Q11").Value = 35608154.91
$ws.Range("R11").Value = 36353781.09
$ws.Range("S11").Value = 37099407.26
$ws.Range("T11").Value = 37845033.44
$ws.Range("U11").Value = 38597538.51
$ws.Range("V11").Value = 39350043.58
$ws.Range("W11").Value = 40102548.66
$ws.Range("X11").Value = 40855053.73
$ws.Range("Y11").Value = 41607558.8
$ws.Range("Z11").Value = 42340492.48
$ws.Range("AA11").Value = 43073426.16
$ws.Range("AB11").Value = 43806359.85
$ws.Range("AC11").Value = 44539293.53
$ws.Range("AD11").Value = 45272227.21
$ws.Range("AE11").Value = 45949578.09
$ws.Range("AF11").Value = 46626928.98
$ws.Range("AG11").Value = 47304279.86
$ws.Range("AH11").Value = 47981630.75
$ws.Range("AI11").Value = 48658981.63
$ws.Range("AJ11").Value = 49257728.04
$ws.Range("AK11").Value = 49856474.46
$ws.Range("AL11").Value = 50455220.87
$ws.Range("AM11").Value = 51053967.28
$ws.Range("AN11").Value = 51652713.7
$ws.Range("AO11").Value = 52154666.72
$ws.Range("AP11").Value = 52656619.75
$ws.Range("AQ11").Value = 53158572.78
$ws.Range("AR11").Value = 53660525.81
$ws.Range("AS11").Value = 54162478.83

# Row 12
$ws.Range("J12").Value = 8432533.607999999
$ws.Range("K12").Value = 8970228.9
$ws.Range("L12").Value = 9549001.752
$ws.Range("M12").Value = 10158400.07
$ws.Range("N12").Value = 10784514.09
$ws.Range("P12").Value = 11987230.9
$ws.Range("Q12").Value = 12630043.9
$ws.Range("R12").Value = 13272856.9
$ws.Range("S12").Value = 13915669.89
$ws.Range("T12").Value = 14558482.89
$ws.Range("U12").Value = 15317480.69
$ws.Range("V12").Value = 16076478.49
$ws.Range("W12").Value = 16835476.28
$ws.Range("X12").Value = 17594474.08
$ws.Range("Y12").Value = 18353471.88
$ws.Range("Z12").Value = 19231050.43
$ws.Range("AA12").Value = 20108628.98
$ws.Range("AB12").Value = 20986207.54
$ws.Range("AC12").Value = 21863786.09
$ws.Range("AD12").Value = 22741364.64
$ws.Range("AE12").Value = 23730211.42
$ws.Range("AF12").Value = 24719058.21
$ws.Range("AG12").Value = 25707904.99
$ws.Range("AH12").Value = 26696751.78
$ws.Range("AI12").Value = 27685598.56
$ws.Range("AJ12").Value = 28778693.21
$ws.Range("AK12").Value = 29871787.86
$ws.Range("AL12").Value = 30964882.51
$ws.Range("AM12").Value = 32057977.16
$ws.Range("AN12").Value = 33151071.81
$ws.Range("AO12").Value = 34339031.9
$ws.Range("AP12").Value = 35526991.98
$ws.Range("AQ12").Value = 36714952.07
$ws.Range("AR12").Value = 37902912.15
$ws.Range("AS12").Value = 39090872.24

Write-Output "Updated rows 2, 10, 11, 12 on strategy_id-0"
